$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.793.00"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "1.881.30"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4725"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3961"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.05"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08058"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("D13").Value = "1.867.02"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.979"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.146"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001050"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "27.807.26"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.541"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "2.102.62"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.39"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.77%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.114"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.594"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9882"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09549"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.452"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.377"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06136"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02258"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.238"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6044"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1908"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.263"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5746"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.950"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.375"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06921"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.05%  "
